$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 (Title): "Testing" " " "custom" " " "properties"
#   -> "Testing " "custom " "properties"
$sh1 = $s.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(1, 8).Text = "Testing "
$tr1.Characters(9, 7).Text = "custom "

# Shape 2 (Subtitle): "This" " " "is" " " "a" " " "subtitle" <br/> <br/> "A." " " "M."
#   -> "This " "is " "a " "subtitle" <br/> <br/> "A. " "M."
$sh2 = $s.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$tr2.Characters(1, 5).Text = "This "
$tr2.Characters(6, 3).Text = "is "
$tr2.Characters(9, 2).Text = "a "
$tr2.Characters(21, 3).Text = "A. "
